$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.028101563453674
$ws.Range("B1").Value = 3.279316663742065
$ws.Range("C1").Value = 3.712397813796997
$ws.Range("D1").Value = 1.993691205978394
$ws.Range("E1").Value = 1.176284432411194
